$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (raw OOXML width = ColumnWidth + 5/6)
$offset = 5/6
$ws.Columns.Item(2).ColumnWidth = 8 - $offset   # column B -> width 8
$ws.Columns.Item(3).ColumnWidth = 8 - $offset   # column C -> width 8
$ws.Columns.Item(7).ColumnWidth = 8 - $offset   # column G -> width 8
$ws.Columns.Item(9).ColumnWidth = 8 - $offset   # column I -> width 8
$ws.Columns.Item(10).ColumnWidth = 8 - $offset   # column J -> width 8
$ws.Columns.Item(11).ColumnWidth = 8 - $offset   # column K -> width 8
$ws.Columns.Item(12).ColumnWidth = 8 - $offset   # column L -> width 8
$ws.Columns.Item(13).ColumnWidth = 8 - $offset   # column M -> width 8
$ws.Columns.Item(15).ColumnWidth = 8 - $offset   # column O -> width 8
$ws.Columns.Item(16).ColumnWidth = 8 - $offset   # column P -> width 8
$ws.Columns.Item(17).ColumnWidth = 8 - $offset   # column Q -> width 8
$ws.Columns.Item(20).ColumnWidth = 9 - $offset   # column T -> width 9
$ws.Columns.Item(22).ColumnWidth = 8 - $offset   # column V -> width 8
$ws.Columns.Item(24).ColumnWidth = 8 - $offset   # column X -> width 8
$ws.Columns.Item(27).ColumnWidth = 8 - $offset   # column AA -> width 8
$ws.Columns.Item(28).ColumnWidth = 8 - $offset   # column AB -> width 8
$ws.Columns.Item(30).ColumnWidth = 8 - $offset   # column AD -> width 8
$ws.Columns.Item(34).ColumnWidth = 8 - $offset   # column AH -> width 8

# Update data values for rows 2-5 (new simulation data)
# Row 2
$ws.Range("A2").Value = 45122.50694444445
$ws.Range("B2").Value = 21.619
$ws.Range("C2").Value = 14.458
$ws.Range("D2").Value = 4.456
$ws.Range("E2").Value = 45.361
$ws.Range("F2").Value = 37.513
$ws.Range("G2").Value = 17.013
$ws.Range("H2").Value = 55.703
$ws.Range("I2").Value = 26.178
$ws.Range("J2").Value = 11.065
$ws.Range("K2").Value = 17.024
$ws.Range("L2").Value = 18.048
$ws.Range("M2").Value = 18.877
$ws.Range("N2").Value = 5.432
$ws.Range("O2").Value = 16.918
$ws.Range("P2").Value = 23.745
$ws.Range("Q2").Value = 14.242
$ws.Range("R2").Value = 3.744
$ws.Range("S2").Value = 2.457
$ws.Range("T2").Value = 250.169
$ws.Range("U2").Value = 47.091
$ws.Range("V2").Value = 15.616
$ws.Range("W2").Value = 31.148
$ws.Range("X2").Value = 16.247
$ws.Range("Y2").Value = 2.407
$ws.Range("Z2").Value = 27.662
$ws.Range("AA2").Value = 13.794
$ws.Range("AB2").Value = 12.753
$ws.Range("AC2").Value = 14.55
$ws.Range("AD2").Value = 18.409
$ws.Range("AE2").Value = 3.64
$ws.Range("AF2").Value = 49.22
$ws.Range("AG2").Value = 8.641
$ws.Range("AH2").Value = 19.523

# Row 3
$ws.Range("A3").Value = 45122.51388888889
$ws.Range("B3").Value = 13.932
$ws.Range("C3").Value = 9.712999999999999
$ws.Range("D3").Value = 1.769
$ws.Range("E3").Value = 29.644
$ws.Range("F3").Value = 24.517
$ws.Range("G3").Value = 10.964
$ws.Range("H3").Value = 43.942
$ws.Range("I3").Value = 16.87
$ws.Range("J3").Value = 7.258
$ws.Range("K3").Value = 10.95
$ws.Range("L3").Value = 11.969
$ws.Range("M3").Value = 12.463
$ws.Range("N3").Value = 3.504
$ws.Range("O3").Value = 10.903
$ws.Range("P3").Value = 15.36
$ws.Range("Q3").Value = 9.381
$ws.Range("R3").Value = 1.593
$ws.Range("S3").Value = 0.997
$ws.Range("T3").Value = 158.641
$ws.Range("U3").Value = 30.563
$ws.Range("V3").Value = 10.064
$ws.Range("W3").Value = 20.221
$ws.Range("X3").Value = 10.801
$ws.Range("Y3").Value = 1.577
$ws.Range("Z3").Value = 20.906
$ws.Range("AA3").Value = 8.888999999999999
$ws.Range("AB3").Value = 8.148
$ws.Range("AC3").Value = 9.443
$ws.Range("AD3").Value = 12.393
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 39.862
$ws.Range("AG3").Value = 5.55
$ws.Range("AH3").Value = 12.582

# Row 4
$ws.Range("A4").Value = 45122.52083333334
$ws.Range("B4").Value = 4.804
$ws.Range("C4").Value = 3.13
$ws.Range("D4").Value = 0.948
$ws.Range("E4").Value = 10.02
$ws.Range("F4").Value = 8.286
$ws.Range("G4").Value = 3.781
$ws.Range("H4").Value = 18.964
$ws.Range("I4").Value = 5.817
$ws.Range("J4").Value = 2.438
$ws.Range("K4").Value = 3.619
$ws.Range("L4").Value = 4.126
$ws.Range("M4").Value = 4.193
$ws.Range("N4").Value = 1.215
$ws.Range("O4").Value = 3.76
$ws.Range("P4").Value = 5.262
$ws.Range("Q4").Value = 3.408
$ws.Range("R4").Value = 0.98
$ws.Range("S4").Value = 0.473
$ws.Range("T4").Value = 49.937
$ws.Range("U4").Value = 10.74
$ws.Range("V4").Value = 3.47
$ws.Range("W4").Value = 6.963
$ws.Range("X4").Value = 3.802
$ws.Range("Y4").Value = 0.613
$ws.Range("Z4").Value = 8.513999999999999
$ws.Range("AA4").Value = 3.065
$ws.Range("AB4").Value = 2.921
$ws.Range("AC4").Value = 3.359
$ws.Range("AD4").Value = 4.221
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 17.392
$ws.Range("AG4").Value = 1.832
$ws.Range("AH4").Value = 4.34

# Row 5
$ws.Range("A5").Value = 45122.52777777778
$ws.Range("B5").Value = 24.02
$ws.Range("C5").Value = 17.69
$ws.Range("D5").Value = 1.37
$ws.Range("E5").Value = 51.93
$ws.Range("F5").Value = 42.94
$ws.Range("G5").Value = 18.9
$ws.Range("H5").Value = 69.05
$ws.Range("I5").Value = 29.09
$ws.Range("J5").Value = 12.9
$ws.Range("K5").Value = 19.29
$ws.Range("L5").Value = 20.95
$ws.Range("M5").Value = 22.01
$ws.Range("N5").Value = 6.04
$ws.Range("O5").Value = 18.8
$ws.Range("P5").Value = 26.74
$ws.Range("Q5").Value = 15.82
$ws.Range("R5").Value = 0.89
$ws.Range("S5").Value = 0.9399999999999999
$ws.Range("T5").Value = 278.8
$ws.Range("U5").Value = 52.47
$ws.Range("V5").Value = 17.35
$ws.Range("W5").Value = 35.28
$ws.Range("X5").Value = 18.73
$ws.Range("Y5").Value = 2.57
$ws.Range("Z5").Value = 34.1
$ws.Range("AA5").Value = 15.33
$ws.Range("AB5").Value = 13.61
$ws.Range("AC5").Value = 15.98
$ws.Range("AD5").Value = 21.94
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 62.44
$ws.Range("AG5").Value = 9.779999999999999
$ws.Range("AH5").Value = 21.69

# Remove row 6 (data now only spans through row 5)
$ws.Rows.Item(6).Delete()

Write-Host "done"